# Apply the crypto price/volume update described by the commit diff.
# Values in column D ("Price") are stored as TEXT in the workbook (inline
# strings), even when they look like plain numbers (e.g. "593.08" or "1.00").
# Excel auto-converts a plain numeric-looking string assigned to .Value into a
# real number (dropping e.g. a trailing ".00"), so for any new Price value that
# looks numeric we force the cell to Text format first, then restore the style
# afterwards so we don't leave stray formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textForcedCells = @(
    "D5", "D6", "D8", "D10", "D11", "D12", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D50", "D51"
)

foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Cell value updates ---
$ws.Range('D2').Value = '67.238.35'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.482.39'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '593.08'
$ws.Range('D6').Value = '177.85'
$ws.Range('E6').Value = '  +3.50%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  +1.36%  '
$ws.Range('D9').Value = '3.484.39'
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').Value = '  +4.75%  '
$ws.Range('D11').Value = '7.07'
$ws.Range('E11').Value = '  -2.26%  '
$ws.Range('D12').Value = '0.434'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '4.083.48'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '31.86'
$ws.Range('E14').Value = '  +9.46%  '
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').Value = '67.251.00'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '0.0000177'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').Value = '3.481.34'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '14.28'
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('D21').Value = '388.30'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').Value = '7.99'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = '73.87'
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '0.536'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = '5.72'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').Value = '0.0000121'
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('D28').Value = '10.34'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('D29').Value = '0.174'
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').Value = '6.12'
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('D32').Value = '1.42'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '23.47'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').Value = '7.35'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '1.60'
$ws.Range('E37').Value = '  -1.61%  '
$ws.Range('D38').Value = '164.47'
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('D39').Value = '0.868'
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  +7.24%  '
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('D43').Value = '4.64'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '2.827.61'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = '26.98'
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').Value = '26.20'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').Value = '0.0722'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').Value = '41.62'
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').Value = '334.50'
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('D51').Value = '1.05'
$ws.Range('E51').Value = '  -2.56%  '

# Restore default styling on the cells we temporarily forced to Text format,
# now that the literal text value is safely stored.
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
